# GetManyCellTypes.xlsx update
# "add cellFormat dateShort, dateLarge, start work on Time"
#
# This commit's net effect on the sheet is a new "monétaire" row (row 23)
# formatted with a new US-dollar style number format, after a brief,
# abandoned experiment with a custom Time format (which is why the
# final currency numFmtId ends up one higher than the previous custom
# format that was added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "start work on Time" ---------------------------------------------
# Trial a custom time number format (work in progress) on the new cell,
# then replace it - this mirrors the numFmtId bookkeeping left behind by
# the real editing session (a Time format was started and abandoned).
$ws.Range("B23").NumberFormat = "[$-F400]h:mm:ss AM/PM;@"

# --- add the new "monétaire" entry -------------------------------------
$ws.Range("A23").Value = "monétaire"
$ws.Range("B23").Value = 91.25
$ws.Range("B23").NumberFormat = "[$$-409]#,##0.00"

# --- update the view: scroll down and select the new cell --------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E22").Select() | Out-Null
